$d = $word.ActiveDocument

# 1) Replace the "UROTAC " heading text with "Test 2".
#    A temporary trailing marker character ("X") is inserted so that we can
#    later anchor the _GoBack bookmark precisely at the end of the new text
#    (immediately after the run, before the paragraph mark) without hitting
#    the runtime's edge-case handling of zero-length ranges placed exactly
#    at a paragraph's text/mark boundary.
$rngReplace = $d.Content
$rngReplace.Find.Execute("UROTAC ", $true, $false, $false, $false, $false, $true, 1, $false, "Test 2X", 2)

# 2) Remove the old _GoBack bookmark (currently sitting right after ##campo##).
If ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks("_GoBack").Delete()
}

# 3) Locate the temporary marker "X" we just inserted.
$rngMarker = $d.Content
$rngMarker.Find.Execute("X", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)

# 4) Add the _GoBack bookmark spanning the marker; once the marker text is
#    removed below the bookmark collapses to a zero-length bookmark sitting
#    right after "Test 2".
$d.Bookmarks.Add("_GoBack", $rngMarker)

# 5) Delete the temporary marker character.
$rngMarker.Text = ""
